$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update team names (column B) for rows 5-10 due to re-ordering of teams
$ws.Range("B5").Value = "Canile Comunale Di Merate"
$ws.Range("B6").Value = "Scroto FC"
$ws.Range("B7").Value = "ReAlcolizzati"
$ws.Range("B8").Value = "FC Tumori"
$ws.Range("B9").Value = "Black Gay United"
$ws.Range("B10").Value = "CSKA LA RISSA"

# Update statistics for rows 5-12 (D..L columns)
$data = @{
    5  = @{ D=11; E=5; F=4; G=2; H=21; I=14; J=7;   K=19; L=771 }
    6  = @{ D=11; E=5; F=3; G=3; H=18; I=17; J=1;   K=18; L=772.5 }
    7  = @{ D=11; E=5; F=2; G=4; H=27; I=15; J=12;  K=17; L=807.5 }
    8  = @{ D=11; E=5; F=2; G=4; H=18; I=22; J=-4;  K=17; L=773.5 }
    9  = @{ D=11; E=4; F=4; G=3; H=13; I=14; J=-1;  K=16; L=743 }
    10 = @{ D=11; E=4; F=3; G=4; H=29; I=22; J=7;   K=15; L=819 }
    11 = @{ D=11; E=4; F=0; G=7; H=13; I=23; J=-10; K=12; L=755 }
    12 = @{ D=11; E=2; F=2; G=7; H=8;  I=20; J=-12; K=8;  L=730 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
